$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the import limit values for 2030/2040/2050 from 0.75 to 0.6
$ws.Range("J5").Value = 0.6
$ws.Range("K5").Value = 0.6
$ws.Range("L5").Value = 0.6

# Add the note about the reset of import limits
$ws.Range("N5").Value = "CGE seems to break if we go to 65% imports."

# Update the active selection to N6 as in the saved file
$ws.Range("N6").Select()
